# Add a new "FAVÖK" (EBITDA) row to both the DATA_RAW sheet (row 26) and
# the "gelir tablosu (çeyreklik)" sheet (row 11), with the same set of
# per-period values across columns B:BU.

$wb = $excel.ActiveWorkbook

$values = @(
    0,0,3101000000,1392000000,0,4725000000,2323000000,2632000000,0,7827000000,
    3979000000,4234000000,0,17340000000,7874000000,9379000000,0,17210000000,
    15384000000,11573000000,0,17463000000,12558000000,9201000000,0,24246000000,
    13593000000,10282000000,0,24171000000,18278000000,14083000000,0,23359000000,
    13987000000,10105000000,0,29019000000,14097000000,10768000000,0,26741000000,
    15894000000,12612000000,0,28385000000,13415000000,11544000000,0,36200000000,
    12853000000,13091000000,0,44185000000,27503000000,24126000000,0,38932000000,
    29979000000,23076000000,0,43221000000,28318000000,22998000000,0,45912000000,
    27900000000,25352000000,0,29589000000,28202000000,0
)

# --- DATA_RAW sheet: new row 26 ---
$wsData = $wb.Worksheets.Item("DATA_RAW")
$wsData.Cells.Item(26, 1).Value = "FAVÖK"
for ($i = 0; $i -lt $values.Length; $i++) {
    $wsData.Cells.Item(26, 2 + $i).Value = $values[$i]
}

# --- "gelir tablosu (çeyreklik)" sheet: new row 11 ---
$wsGelir = $wb.Worksheets.Item("gelir tablosu (çeyreklik)")
$wsGelir.Cells.Item(11, 1).Value = "FAVÖK"
for ($i = 0; $i -lt $values.Length; $i++) {
    $wsGelir.Cells.Item(11, 2 + $i).Value = $values[$i]
}
